# heliaphen_experiments.xlsx
# - widen the data columns (A..I) to make room for the new download
#   link / water-deficit slider controls that were added to the sheet
# - bump the "weight_dead" (I) values for the water-deficit slider rows
#   (I3:I8) from 450 -> 1300; I2 is left untouched
# - move the active selection from D3 to I3 (new slider column)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- widen columns A:I -----------------------------------------------
# ColumnWidth is expressed in "characters" and gets rounded by the COM
# layer to the nearest 1/6th before the 5/6-character padding is added
# back in on save, so we pick the nearest representable input for each
# target width.
$ws.Columns.Item(1).ColumnWidth = 14.5              # A
$ws.Columns.Item(2).ColumnWidth = 161.5             # B
$ws.Range("C1:D1").EntireColumn.ColumnWidth = 11.6666666666667   # C:D
$ws.Range("E1:F1").EntireColumn.ColumnWidth = 12.3333333333333   # E:F
$ws.Range("G1:H1").EntireColumn.ColumnWidth = 5.16666666666667   # G:H
$ws.Columns.Item(9).ColumnWidth = 16.5              # I

# --- weight_dead slider values (column I, rows 3-8) -------------------
$ws.Range("I3:I8").Value = 1300

# --- selection moves to the new slider cell ---------------------------
$ws.Range("I3").Select()
